# Append: 2026-01-30 18:50 JST
# Update the "取得日時" (acquired timestamp) column (A) for all data rows
# on the "ランサーズ" sheet from the old timestamp to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2026-01-30 18:38:56"
$newTimestamp = "2026-01-30 18:50:16"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
